$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11, columns B-G.
# Each row's B:F values "shift down" from the previous row (a new forecast
# point was prepended and the oldest one dropped), and the N count (col G)
# increases by one to reflect the additional observation.
$data = @(
    @(0.1184146901506048, 0.5185973060008381, 0.4691159408045538, 0.6849203901217672, 0.698284046247675, 15),
    @(0.1872890445127373, 0.5754152234511226, 0.6110775347743189, 0.7817144841784108, 0.7875963862927191, 14),
    @(0.1843166320775897, 0.6273081474520239, 0.5779702429309931, 0.7602435418541831, 0.7676788112505872, 13),
    @(0.3221619536241352, 0.5655998917628668, 0.5382575468002283, 0.7336603756509059, 0.6884521947652315, 12),
    @(0.2262967854219969, 0.4712852108734478, 0.2956388450555676, 0.5437268110508876, 0.5185281775945857, 11),
    @(0.3318332233303713, 0.399431789294712,  0.2898093932768167, 0.53833947772462,   0.4468359196381605, 10),
    @(0.3032448163457492, 0.4442789123723634, 0.3348453944406773, 0.5786582708651776, 0.5227322189947009, 9),
    @(0.3395226469355723, 0.4954646948149822, 0.3955394669362153, 0.6289192849135852, 0.5659518799878692, 8),
    @(0.2802560178301312, 0.3574868270362884, 0.1586788892857869, 0.3983451886062977, 0.3057635732311509, 7),
    @(0.3351585394422735, 0.4364782826832794, 0.3756151435636499, 0.612874492505317,  0.5620860044540215, 6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
